# Daily attendance processing - 2025-12-28 16:36:00
# Normalizes the "Recorded By" (column G) audit-trail lists so that the
# most-recent recorder (which was appended last) is promoted to the front
# of the comma-separated list, unless the list already ends with "System".

# Case-sensitive string equality check. The emulated -ceq/-cne operators in
# this runtime behave case-insensitively, so compare character codes directly
# to reliably distinguish "System" from "system".
function Test-ExactEquals($s1, $s2) {
    if ($s1.Length -ne $s2.Length) { return $false }
    for ($i = 0; $i -lt $s1.Length; $i++) {
        if ([int][char]$s1[$i] -ne [int][char]$s2[$i]) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 157; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ([string]::IsNullOrEmpty($value)) { continue }

    $parts = $value -split ",\s*"

    if ($parts.Count -gt 1 -and -not (Test-ExactEquals $parts[$parts.Count - 1] "System")) {
        $lastPart = $parts[$parts.Count - 1]
        $rest = $parts[0..($parts.Count - 2)]
        $newParts = @($lastPart) + $rest
        $newValue = [string]::Join(", ", $newParts)
        $cell.Value2 = $newValue
    }
}
